# Update Data.xlsx: replace "Capsule Image" (column C) hyperlink URLs with
# locally-referenced image paths, drop the hyperlink on column C (it is no
# longer a link target, just a plain filename), and keep the "Sleeve Image"
# (column D) hyperlinks intact for rows 2-37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Replace Capsule Image (column C) cell values with local image paths ---
$ws.Cells.Item(2, 3).Value = 'Images/OL1_Ispirazione Napoli_Original_Capsule.png'
$ws.Cells.Item(3, 3).Value = 'Images/OL2_Kazaar_Original_Capsule.png'
$ws.Cells.Item(4, 3).Value = 'Images/OL3_Ristretto_Original_Capsule.png'
$ws.Cells.Item(5, 3).Value = 'Images/OL4_Ristretto Decaffeinato_Original_Capsule.png'
$ws.Cells.Item(6, 3).Value = 'Images/OL5_Arpeggio_Original_Capsule.png'
$ws.Cells.Item(7, 3).Value = 'Images/OL6_Arpeggio Decaffeinato_Original_Capsule.png'
$ws.Cells.Item(8, 3).Value = 'Images/OL7_Inspirazione Venezia_Original_Capsule.png'
$ws.Cells.Item(9, 3).Value = 'Images/OL8_Inspirazione Roma_Original_Capsule.png'
$ws.Cells.Item(10, 3).Value = 'Images/OL9_Livanto_Original_Capsule.png'
$ws.Cells.Item(11, 3).Value = 'Images/OL10_Cape Town Lungo_Original_Capsule.png'
$ws.Cells.Item(12, 3).Value = 'Images/OL11_Miami Espresso_Original_Capsule.png'
$ws.Cells.Item(13, 3).Value = 'Images/OL12_Rio De Janeiro Espresso_Original_Capsule.png'
$ws.Cells.Item(14, 3).Value = 'Images/OL13_Istanbul Espresso_Original_Capsule.png'
$ws.Cells.Item(15, 3).Value = 'Images/OL14_Stockholm Lungo_Original_Capsule.png'
$ws.Cells.Item(16, 3).Value = 'Images/OL15_Paris Espresso_Original_Capsule.png'
$ws.Cells.Item(17, 3).Value = 'Images/OL16_Vienna Lungo_Original_Capsule.png'
$ws.Cells.Item(18, 3).Value = 'Images/OL17_Tokyo Lungo_Original_Capsule.png'
$ws.Cells.Item(19, 3).Value = 'Images/OL18_Shanghai Lungo_Original_Capsule.png'
$ws.Cells.Item(20, 3).Value = 'Images/OL19_Buenos Aires Lungo_Original_Capsule.png'
$ws.Cells.Item(21, 3).Value = 'Images/OL20_India_Original_Capsule.png'
$ws.Cells.Item(22, 3).Value = 'Images/OL21_Indonesia - Fairtrade_Original_Capsule.png'
$ws.Cells.Item(23, 3).Value = 'Images/OL22_Colombia_Original_Capsule.png'
$ws.Cells.Item(24, 3).Value = 'Images/OL23_Peru Organic_Original_Capsule.png'
$ws.Cells.Item(25, 3).Value = 'Images/OL24_Nicaragua_Original_Capsule.png'
$ws.Cells.Item(26, 3).Value = 'Images/OL25_Ethiopia_Original_Capsule.png'
$ws.Cells.Item(27, 3).Value = 'Images/OL26_Cioccolatino_Original_Capsule.png'
$ws.Cells.Item(28, 3).Value = 'Images/OL27_Vaniglia_Original_Capsule.png'
$ws.Cells.Item(29, 3).Value = 'Images/OL28_Nocciola_Original_Capsule.png'
$ws.Cells.Item(30, 3).Value = 'Images/OL29_Caramello_Original_Capsule.png'
$ws.Cells.Item(31, 3).Value = 'Images/OL30_Corto_Original_Capsule.png'
$ws.Cells.Item(32, 3).Value = 'Images/OL31_Scuro_Original_Capsule.png'
$ws.Cells.Item(33, 3).Value = 'Images/OL32_Chiaro_Original_Capsule.png'
$ws.Cells.Item(34, 3).Value = 'Images/OL33_Capriccio_Original_Capsule.png'
$ws.Cells.Item(35, 3).Value = 'Images/OL34_Cosi_Original_Capsule.png'
$ws.Cells.Item(36, 3).Value = 'Images/OL35_Volluto_Original_Capsule.png'
$ws.Cells.Item(37, 3).Value = 'Images/OL36_Volluto Decaffeinato_Original_Capsule.png'

# --- Step 2: Drop every hyperlink on the sheet. (This engine only supports bulk
#     removal via the sheet-level Hyperlinks collection; per-cell deletion is a
#     no-op here.) We recreate the ones we want to keep (column D) immediately after. ---
$ws.Hyperlinks.Delete()

# --- Step 3: Column C is plain text now, so reset its style back to Normal ---
$ws.Cells.Item(2, 3).Style = "Normal"
$ws.Cells.Item(3, 3).Style = "Normal"
$ws.Cells.Item(4, 3).Style = "Normal"
$ws.Cells.Item(5, 3).Style = "Normal"
$ws.Cells.Item(6, 3).Style = "Normal"
$ws.Cells.Item(7, 3).Style = "Normal"
$ws.Cells.Item(8, 3).Style = "Normal"
$ws.Cells.Item(9, 3).Style = "Normal"
$ws.Cells.Item(10, 3).Style = "Normal"
$ws.Cells.Item(11, 3).Style = "Normal"
$ws.Cells.Item(12, 3).Style = "Normal"
$ws.Cells.Item(13, 3).Style = "Normal"
$ws.Cells.Item(14, 3).Style = "Normal"
$ws.Cells.Item(15, 3).Style = "Normal"
$ws.Cells.Item(16, 3).Style = "Normal"
$ws.Cells.Item(17, 3).Style = "Normal"
$ws.Cells.Item(18, 3).Style = "Normal"
$ws.Cells.Item(19, 3).Style = "Normal"
$ws.Cells.Item(20, 3).Style = "Normal"
$ws.Cells.Item(21, 3).Style = "Normal"
$ws.Cells.Item(22, 3).Style = "Normal"
$ws.Cells.Item(23, 3).Style = "Normal"
$ws.Cells.Item(24, 3).Style = "Normal"
$ws.Cells.Item(25, 3).Style = "Normal"
$ws.Cells.Item(26, 3).Style = "Normal"
$ws.Cells.Item(27, 3).Style = "Normal"
$ws.Cells.Item(28, 3).Style = "Normal"
$ws.Cells.Item(29, 3).Style = "Normal"
$ws.Cells.Item(30, 3).Style = "Normal"
$ws.Cells.Item(31, 3).Style = "Normal"
$ws.Cells.Item(32, 3).Style = "Normal"
$ws.Cells.Item(33, 3).Style = "Normal"
$ws.Cells.Item(34, 3).Style = "Normal"
$ws.Cells.Item(35, 3).Style = "Normal"
$ws.Cells.Item(36, 3).Style = "Normal"
$ws.Cells.Item(37, 3).Style = "Normal"

# --- Step 4: Recreate the Sleeve Image (column D) hyperlinks, preserving targets and order ---
$ws.Hyperlinks.Add($ws.Cells.Item(2, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/ispirazione-napoli_L.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3, 4), 'https://m.media-amazon.com/images/I/51vv3BK5STL._AC_UF1000,1000_QL80_.jpg') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(4, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/ispirazione-ristretto_L.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(5, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/ispirazione-ristretto-decaffeinato_S.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(6, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/ispirazione-firenze-arpeggio_L.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(7, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/ispirazione-firenze-arpeggio-decaffeinato_L.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(8, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/ispirazione-venezia_S.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(9, 4), 'https://www.nespresso.si/files/thumbs/files/images/product/thumbs_600/7747-80_2_600_600px.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(10, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/ispirazione-genova-livanto_XL.png?impolicy=medium&imwidth=824&imdensity=1') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(11, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/cape-town-envivo-lungo_XL.png?impolicy=medium&imwidth=824&imdensity=1') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(12, 4), 'https://www.nespresso.com/ecom/medias/sys_master/public/15143389757470/AT-B2C-2022-WEX-Miami-Image-Set-Mobile.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(13, 4), 'https://www.nespresso.com/shared_res/agility/n-components/B2C-enriched-pdp-wex-2022/rio-de-janeiro-espresso/main-image_L.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(14, 4), 'https://www.nespresso.com/shared_res/agility/n-components/B2C-enriched-pdp-wex-2022/istanbul-espresso/main-image_S.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(15, 4), 'https://www.nespresso.com/ecom/medias/sys_master/public/15941486116894/ol-coffee-sleeves-stockholm-16-9.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(16, 4), 'https://www.nespresso.com/shared_res/agility/n-components/B2C-enriched-pdp-wex-2022/paris-espresso/main-image_S.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(17, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/vienna-linizio-lungo_L.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(18, 4), 'https://i.ebayimg.com/images/g/cTQAAOSwjrpkAaJ9/s-l1200.webp') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(19, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/shangai-lungo_S.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(20, 4), 'https://www.nespresso.si/files/thumbs/files/images/product/thumbs_600/7757-80_2_600_600px.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(21, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/india_L.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(22, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/indonesia_S.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(23, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/colombia_L.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(24, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/peru_S.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(25, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/nicaragua_XL.png?impolicy=medium&imwidth=824&imdensity=1') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(26, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/ethiopia_L.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(27, 4), 'https://www.nespresso.si/files/thumbs/files/images/product/thumbs_600/cioccolatino-nespresso-capsules-sleeve_600_600px.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(28, 4), 'https://www.nespresso.com/shared_res/agility/global/coffees/ol/sku-main-info-product/vaniglia_2x.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(29, 4), 'https://www.nespresso.com/shared_res/agility/global/coffees/ol/sku-main-info-product/nocciola_2x.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(30, 4), 'https://www.nespresso.com/shared_res/agility/global/coffees/ol/sku-main-info-product/caramello_2x.png?impolicy=medium&imwidth=824&imdensity=1') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(31, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/corto_L.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(32, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/scuro_XL.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(33, 4), 'https://www.nespresso.com/shared_res/agility/n-components/pdp/sku-main-info/coffee-sleeves/ol/chiaro_L.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(34, 4), 'https://www.nespresso.com/shared_res/agility/commons/img/coffees/OL/composition/ol_coffee-sleeves_capriccio_16-9_2x.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(35, 4), 'https://www.nespresso.com/shared_res/agility/global/coffees/ol/sku-main-info-product/cosi_2x.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(36, 4), 'https://www.nespresso.com/shared_res/agility/commons/img/coffees/OL/composition/ol_coffee-sleeves_volluto_2x.png') | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(37, 4), 'https://down-ph.img.susercontent.com/file/384c222616f53cf0a3676311b8b8a33d') | Out-Null

# --- Step 5: Re-apply the Hyperlink style to column D (Add() above resets it to an
#     equivalent but separate style record; restoring the named style keeps it tidy) ---
$ws.Cells.Item(2, 4).Style = "Hyperlink"
$ws.Cells.Item(3, 4).Style = "Hyperlink"
$ws.Cells.Item(4, 4).Style = "Hyperlink"
$ws.Cells.Item(5, 4).Style = "Hyperlink"
$ws.Cells.Item(6, 4).Style = "Hyperlink"
$ws.Cells.Item(7, 4).Style = "Hyperlink"
$ws.Cells.Item(8, 4).Style = "Hyperlink"
$ws.Cells.Item(9, 4).Style = "Hyperlink"
$ws.Cells.Item(10, 4).Style = "Hyperlink"
$ws.Cells.Item(11, 4).Style = "Hyperlink"
$ws.Cells.Item(12, 4).Style = "Hyperlink"
$ws.Cells.Item(13, 4).Style = "Hyperlink"
$ws.Cells.Item(14, 4).Style = "Hyperlink"
$ws.Cells.Item(15, 4).Style = "Hyperlink"
$ws.Cells.Item(16, 4).Style = "Hyperlink"
$ws.Cells.Item(17, 4).Style = "Hyperlink"
$ws.Cells.Item(18, 4).Style = "Hyperlink"
$ws.Cells.Item(19, 4).Style = "Hyperlink"
$ws.Cells.Item(20, 4).Style = "Hyperlink"
$ws.Cells.Item(21, 4).Style = "Hyperlink"
$ws.Cells.Item(22, 4).Style = "Hyperlink"
$ws.Cells.Item(23, 4).Style = "Hyperlink"
$ws.Cells.Item(24, 4).Style = "Hyperlink"
$ws.Cells.Item(25, 4).Style = "Hyperlink"
$ws.Cells.Item(26, 4).Style = "Hyperlink"
$ws.Cells.Item(27, 4).Style = "Hyperlink"
$ws.Cells.Item(28, 4).Style = "Hyperlink"
$ws.Cells.Item(29, 4).Style = "Hyperlink"
$ws.Cells.Item(30, 4).Style = "Hyperlink"
$ws.Cells.Item(31, 4).Style = "Hyperlink"
$ws.Cells.Item(32, 4).Style = "Hyperlink"
$ws.Cells.Item(33, 4).Style = "Hyperlink"
$ws.Cells.Item(34, 4).Style = "Hyperlink"
$ws.Cells.Item(35, 4).Style = "Hyperlink"
$ws.Cells.Item(36, 4).Style = "Hyperlink"
$ws.Cells.Item(37, 4).Style = "Hyperlink"

